$wb = $excel.ActiveWorkbook

# ALC row 42
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 419.8889
$ws.Range("I42").Value = 576
$ws.Range("J42").Value = 224.75
$ws.Range("K42").Value = 1728
$ws.Range("L42").Value = 674.25
$ws.Range("M42").Value = -1498
$ws.Range("N42").Value = -1134.25

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5953.048
$ws.Range("I74").Value = 4860.5
$ws.Range("K74").Value = 4860.5
$ws.Range("M74").Value = -3924.5

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3065.4
$ws.Range("I76").Value = 2961.6223
$ws.Range("K76").Value = 2961.6223
$ws.Range("M76").Value = -2646.6223

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5953.048
$ws.Range("I77").Value = 4860.5
$ws.Range("K77").Value = 24302.5
$ws.Range("M77").Value = -19622.5

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3065.4
$ws.Range("I79").Value = 2961.6223
$ws.Range("K79").Value = 2961.6223
$ws.Range("M79").Value = -1869.6223

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1374.2113
$ws.Range("J112").Value = 1386.7
$ws.Range("L112").Value = 4160.1
$ws.Range("N112").Value = -6376.1

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 22818620
$ws.Range("I132").Value = 25742752
$ws.Range("K132").Value = 77228256
$ws.Range("M132").Value = -77225726

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 975180
$ws.Range("I137").Value = 2805677.2
$ws.Range("J137").Value = 2728.3125
$ws.Range("K137").Value = 8417031.600000001
$ws.Range("L137").Value = 8184.9375
$ws.Range("M137").Value = -8414481.600000001
$ws.Range("N137").Value = -13284.9375

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1571
$ws.Range("I138").Value = 1209.5294
$ws.Range("J138").Value = 2800
$ws.Range("K138").Value = 3628.5882
$ws.Range("L138").Value = 8400
$ws.Range("M138").Value = 1511.4118
$ws.Range("N138").Value = -18680

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 8660664
$ws.Range("I63").Value = 15392314
$ws.Range("J63").Value = 5684
$ws.Range("K63").Value = 15392314
$ws.Range("L63").Value = 5684
$ws.Range("M63").Value = -15391628
$ws.Range("N63").Value = -7056

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 8660664
$ws.Range("I66").Value = 15392314
$ws.Range("J66").Value = 5684
$ws.Range("K66").Value = 76961570
$ws.Range("L66").Value = 28420
$ws.Range("M66").Value = -76958138
$ws.Range("N66").Value = -35284

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 248938.52
$ws.Range("I74").Value = 410495.88
$ws.Range("J74").Value = 1850.7646
$ws.Range("K74").Value = 410495.88
$ws.Range("L74").Value = 1850.7646
$ws.Range("M74").Value = -409621.88
$ws.Range("N74").Value = -3598.7646

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 248938.52
$ws.Range("I77").Value = 410495.88
$ws.Range("J77").Value = 1850.7646
$ws.Range("K77").Value = 2052479.4
$ws.Range("L77").Value = 9253.823
$ws.Range("M77").Value = -2048111.4
$ws.Range("N77").Value = -17989.823

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1786.6329
$ws.Range("I105").Value = 1670.5286
$ws.Range("K105").Value = 1670.5286
$ws.Range("M105").Value = 76.4713999999999

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3400.7856
$ws.Range("I134").Value = 1742.8572
$ws.Range("J134").Value = 5058.7144
$ws.Range("K134").Value = 5228.571599999999
$ws.Range("L134").Value = 15176.1432
$ws.Range("M134").Value = -2693.571599999999
$ws.Range("N134").Value = -20246.1432

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1298.2941
$ws.Range("I94").Value = 649.2222
$ws.Range("K94").Value = 649.2222
$ws.Range("M94").Value = -198.2222

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1522.4375
$ws.Range("I134").Value = 994.5417
$ws.Range("J134").Value = 3106.125
$ws.Range("K134").Value = 2983.6251
$ws.Range("L134").Value = 9318.375
$ws.Range("M134").Value = -448.6251000000002
$ws.Range("N134").Value = -14388.375

# CUL row 97
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 870.53845
$ws.Range("I97").Value = 803
$ws.Range("J97").Value = 1022.5
$ws.Range("K97").Value = 2409
$ws.Range("L97").Value = 3067.5
$ws.Range("M97").Value = -1913
$ws.Range("N97").Value = -4059.5

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 578.8182
$ws.Range("I98").Value = 474.6
$ws.Range("J98").Value = 665.6667
$ws.Range("K98").Value = 1423.8
$ws.Range("L98").Value = 1997.0001
$ws.Range("M98").Value = 74.19999999999982
$ws.Range("N98").Value = -4993.0001

# CUL row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 4928
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 4928
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 14784
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -19276

# CUL row 107
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 9824809
$ws.Range("I107").Value = 356.13635
$ws.Range("J107").Value = 17277842
$ws.Range("K107").Value = 1068.40905
$ws.Range("L107").Value = 51833526
$ws.Range("M107").Value = 851.59095
$ws.Range("N107").Value = -51837366

# CUL row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 7928.357
$ws.Range("I109").Value = 2999.6667
$ws.Range("J109").Value = 9272.546
$ws.Range("K109").Value = 8999.000100000001
$ws.Range("L109").Value = 27817.638
$ws.Range("M109").Value = -7959.000100000001
$ws.Range("N109").Value = -29897.638

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1612.3658
$ws.Range("I121").Value = 509
$ws.Range("J121").Value = 1820.2463
$ws.Range("K121").Value = 1527
$ws.Range("L121").Value = 5460.7389
$ws.Range("M121").Value = -217
$ws.Range("N121").Value = -8080.7389

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 2395.689
$ws.Range("I122").Value = 672.0625
$ws.Range("J122").Value = 3346.6553
$ws.Range("K122").Value = 6048.5625
$ws.Range("L122").Value = 30119.8977
$ws.Range("M122").Value = -3598.5625
$ws.Range("N122").Value = -35019.8977

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6057.463
$ws.Range("I70").Value = 5557.8823
$ws.Range("K70").Value = 5557.8823
$ws.Range("M70").Value = -5287.8823

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6057.463
$ws.Range("I73").Value = 5557.8823
$ws.Range("K73").Value = 5557.8823
$ws.Range("M73").Value = -4621.8823

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 35716812
$ws.Range("I80").Value = 125001250
$ws.Range("J80").Value = 3037.2
$ws.Range("K80").Value = 125001250
$ws.Range("L80").Value = 3037.2
$ws.Range("M80").Value = -125000252
$ws.Range("N80").Value = -5033.2

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 35716812
$ws.Range("I83").Value = 125001250
$ws.Range("J83").Value = 3037.2
$ws.Range("K83").Value = 625006250
$ws.Range("L83").Value = 15186
$ws.Range("M83").Value = -625001258
$ws.Range("N83").Value = -25170
